$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Sponsored
Launching 2 BHK in Thane West - Homes in Pokhran Rd by...
raymondtenxera.com
https://www.raymondtenxera.com › official-site › brand
New Launch Homes by Raymond with 38 Habitable Floors, 26500 SqFt Clubhouse, 40+ Amenities. Experience a futuristic lifestyle with Raymond Realty''s Spacious 2 BHK homes...'
$ws.Range("B2").Value = 'Raymond Limited'
$ws.Range("C2").Value = 'India'

$ws.Range("A3").Value = 'Sponsored
Visit Address By GS Thane - Consult an expert & visit site
raymonds-addressbygs.com
https://www.raymonds-addressbygs.com
Bookings Open Addres By GS Thane 6.1 Acre, 2/3/4 Bhk 1.30 Cr Ask Expert & Visit Site. Booking Open For Limited Time Only, Sign Up & Get Instant Call Back Now. Early Buy Discount.
View Pricing · Price ₹ / BHK/ Area · View The Gallery · Site & Layout Plan · Grand Amenities
2 Bhk - 615 Sq. Ft. - ₹1.61 - Price In Cr
 · 
More'
$ws.Range("B3").Value = 'Home Bazaar Services Pvt Ltd'
$ws.Range("C3").Value = 'India'

$ws.Range("A4").Value = 'Sponsored
La Vie at Uptown Urbania | 2&3 BHK at ₹1.49 Cr (All Incl)
Rustomjee La-Vie
https://www.rustomjee-lavie.com
Book 2 & 3 BHK homes from ₹1.49 Cr (All Incl) at Rustomjee La Vie, Thane (W) 2&3BHK from ₹1.49 Cr (All Incl) at Rustomjee La Vie. Pay 20% Now & Nothing Till Jan''25. Luxury flats. Leisure zones.'
$ws.Range("B4").Value = 'Kapstone Construction Pvt Ltd'
$ws.Range("C4").Value = 'India'

$ws.Range("A5").Value = 'Sponsored
The Address By GS Thane | 3, 4 & 4.5 BHK ₹2.59Cr*
theaddressbygs-thane.in
https://www.theaddressbygs-thane.in
Launching The Address by GS at Pokhran Road Thane. Price Starts at ₹2.59 Cr*. Book Now.'
$ws.Range("B5").Value = 'Prop Solutions4u'
$ws.Range("C5").Value = 'India'

$ws.Range("A6").Value = 'Sponsored
TenX Habitat Thane Launch | 2, 3 & 4 BHK Starts @ 1.41Cr*
homesfy-property.co.in
https://www.homesfy-property.co.in › tenx › thane
326,500 sq ft clubhouse | 2 and 3 BHK residences near Viviana Mall @ 1.41 CR* Launching Tenx Habitat at Thane. Lavish 2, 3 & 4 BHK Starting At 1.41Cr* Highlights: Chat Option Available, Floor Plan Available, Brochure Available.
Price Plan · Our Price · Browse Prices · Floor Plans · Floor Plan · View Gallery'
$ws.Range("B6").Value = 'Homesfy Realty Limited'
$ws.Range("C6").Value = 'India'

$ws.Range("A7").Value = 'Sponsored
Launching Final Tower Asteria
luxeoffplans.com
https://www.luxeoffplans.com
2, 3 & 4 BHK Apartments — Spacious 2, 3 & 4 BHK Homes at Thane. Sample Flat Ready. Download Brochure. Explore Project Details, Speak To Site Expert, Know More.'
$ws.Range("B7").Value = 'RABS NET SOLUTIONS PVT LTD'
$ws.Range("C7").Value = 'India'

$ws.Range("A8").Value = 'Sponsored
Exclusive Offers On Site
raymondsrealtythane.com
https://www.raymondsrealtythane.com
Book Online Presentation Today — Luxurious 2, 3, 4 & 4.5 BHK Starts At ₹ 1.30 Cr All Inc | Flexi Payment Plan Available. The Address By GS 2.0 Offers Luxurious 2, 3, 4 & 4.5 BHK Home With Balcony. Pay Just 20% Now. Book a free Site Visit. Easy Payment Plan. Amenities: Senior Citizen Area, Yoga Path.'
$ws.Range("B8").Value = 'zuber khan'
$ws.Range("C8").Value = 'India'

$ws.Range("A9").Value = 'Sponsored
The Address By GS Thane - 3, 4 & 4.5 BHK ₹2.59Cr*
theaddressbygs-thane.in
https://www.theaddressbygs-thane.in
Launching The Address by GS at Pokhran Road Thane. Price Starts at ₹2.59 Cr*. Book Now.
Call us'
$ws.Range("B9").Value = 'Prop Solutions4u'
$ws.Range("C9").Value = 'India'

$ws.Range("A10").Value = 'Sponsored
Visit Address Tower Thane - Consult an expert & visit site
raymonds-addressbygs.com
https://www.raymonds-addressbygs.com
Bookings Open Addres Tower Thane 6.1 Acre, 2/3/4 Bhk 1.30 Cr Ask Expert & Visit Site'
$ws.Range("B10").Value = 'Home Bazaar Services Pvt Ltd'
$ws.Range("C10").Value = 'India'

$ws.Range("A11").Value = 'Sponsored
TenX Habitat Thane Launch - Zero Stamp Duty Offer
homesfy-property.co.in
https://www.homesfy-property.co.in › tenx › thane
326,500 sq ft clubhouse | 2 and 3 BHK residences near Viviana Mall @ 1.41 CR*'
$ws.Range("B11").Value = 'Homesfy Realty Limited'
$ws.Range("C11").Value = 'India'

$ws.Range("A12").Value = 'Sponsored
La Vie at Uptown Urbania - 2&3 BHK at ₹1.49 Cr (All Incl)
Rustomjee La-Vie
https://www.rustomjee-lavie.com
Book 2 & 3 BHK homes from ₹1.49 Cr (All Incl) at Rustomjee La Vie, Thane (W...'
$ws.Range("B12").Value = 'Kapstone Construction Pvt Ltd'
$ws.Range("C12").Value = 'India'

$ws.Range("A13").Value = 'Sponsored
TenX Habitat Thane - Zero Stamp Duty Offer
homesfy-property.co.in
https://www.homesfy-property.co.in › tenx › thane
326,500 sq ft clubhouse | 2 and 3 BHK residences near Viviana Mall @ 1.41 CR*'
$ws.Range("B13").Value = 'Homesfy Realty Limited'
$ws.Range("C13").Value = 'India'

$ws.Range("A14").Value = 'Sponsored
The Address By GS Pokhran Road | 3, 4, 4.5 BHK Price ₹2.59Cr*
theaddressbygs-thane.in
https://www.theaddressbygs-thane.in
Launching The Address by GS at Pokhran Road Thane. Price Starts at ₹2.59 Cr*. Book Now.'
$ws.Range("B14").Value = 'Prop Solutions4u'
$ws.Range("C14").Value = 'India'

$ws.Range("A15").Value = 'Sponsored
Visit Address Tower Thane | Get expert advice & visit site
raymonds-addressbygs.com
https://www.raymonds-addressbygs.com
Bookings Open Addres Tower Thane 6.1 Acre, 2/3/4 Bhk 1.30 Cr Ask Expert & Visit Site'
$ws.Range("B15").Value = 'Home Bazaar Services Pvt Ltd'
$ws.Range("C15").Value = 'India'

$ws.Range("A16").Value = 'Sponsored
Thane''s Finest Homes @Thane | Ready World Class Amenities
raymondsrealtythane.com
https://www.raymondsrealtythane.com
Luxurious 2, 3, 4 & 4.5 BHK Starts At ₹ 1.30 Cr All Inc | Flexi Payment Plan Available. The...'
$ws.Range("B16").Value = 'zuber khan'
$ws.Range("C16").Value = 'India'

$ws.Range("A17").Value = 'Sponsored
2 BHK with Balcony at Raymond - Presenting Homes in Thane...
raymondtenxera.com
https://www.raymondtenxera.com › official-site › brand
New Launch Homes by Raymond with 38 Habitable Floors, 26500 SqFt Clubhouse, 40+ Amenities'
$ws.Range("B17").Value = 'Raymond Limited'
$ws.Range("C17").Value = 'India'

$ws.Range("A18").Value = 'Sponsored
TenX Habitat Thane | Zero Stamp Duty Offer
homesfy-property.co.in
https://www.homesfy-property.co.in › tenx › thane
326,500 sq ft clubhouse | 2 and 3 BHK residences near Viviana Mall @ 1.41 CR*'
$ws.Range("B18").Value = 'Homesfy Realty Limited'
$ws.Range("C18").Value = 'India'

$ws.Range("A19").Value = 'Sponsored
The Address By GS Thane | 3, 4, 4.5 BHK Price ₹2.59Cr*
theaddressbygs-thane.in
https://www.theaddressbygs-thane.in
Launching The Address by GS at Pokhran Road Thane. Price Starts at ₹2.59 Cr*. Book Now.'
$ws.Range("B19").Value = 'Prop Solutions4u'
$ws.Range("C19").Value = 'India'

$ws.Range("A20").Value = 'Sponsored
RTMI Flats with Zero GST* - Grand Clubhouse @Raymond TenX
tenxhabitatraymondrealty.com
https://www.tenxhabitatraymondrealty.com
4,200 Sq.Ft. Multipurpose Hall | 28 Seater Mini Theatre | 2,400 Sq.Ft. Fully Equipped Gym'
$ws.Range("B20").Value = 'Raymond Limited'
$ws.Range("C20").Value = 'India'
